$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held a duplicate "5983729 - Fernando Vernilli Junior"
# value (B13/C13, with no label in column A) is removed entirely; Excel
# shifts every following row up by one and keeps everybody else's row
# heights / styles intact.
$ws.Rows.Item(13).Delete()

# After the shift, fix up the handful of description cells whose text
# content changed (the column-A labels already line up correctly).
$ws.Range("B10").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C10").Value = "5983729 - Fernando Vernilli Junior"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# B15/C15 need the literal text "01/01/2011" (same string already used in
# B8/C8). Assigning that string directly via .Value gets auto-converted to
# a date serial by the host, so instead copy the already-correct text
# value over from B8/C8, after first copying the destination's own
# pre-existing cell format back onto itself so paste-values can't drag in
# a foreign style.
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (no-op, keeps s=2)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("C15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("B18").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C18").Value = "5983729 - Fernando Vernilli Junior"

$ws.Range("B19").Value = "A avaliação será constituída por aulas expositivas, aulas de exercícios e laboratórios. Serão aplicadas pelo menos duas provas."
$ws.Range("C19").Value = "A avaliação será constituída por aulas expositivas, aulas de exercícios e laboratórios. Serão aplicadas pelo menos duas provas."

$ws.Range("B20").Value = "A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final)."
$ws.Range("C20").Value = "A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final)."

$ws.Range("B21").Value = "A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2."
$ws.Range("C21").Value = "A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2."
